$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add D5 with the same value as B5 ("30min")
$ws.Range("D5").Value = "30min"

# Update the selection to D6 (matches the diff's sheetView selection change)
$ws.Range("D6").Select()
